$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3893.3
$ws.Range("J64").Value = 3993.2666
$ws.Range("L64").Value = 3993.2666
$ws.Range("N64").Value = -4489.2666
$ws.Range("H67").Value = 3893.3
$ws.Range("J67").Value = 3993.2666
$ws.Range("L67").Value = 3993.2666
$ws.Range("N67").Value = -5709.2666
$ws.Range("H112").Value = 1348.2273
$ws.Range("I112").Value = 599.75
$ws.Range("J112").Value = 1514.5555
$ws.Range("K112").Value = 1799.25
$ws.Range("L112").Value = 4543.666499999999
$ws.Range("M112").Value = -691.25
$ws.Range("N112").Value = -6759.666499999999
$ws.Range("H132").Value = 2360.6516
$ws.Range("I132").Value = 2055.0557
$ws.Range("J132").Value = 3735.8333
$ws.Range("K132").Value = 6165.1671
$ws.Range("L132").Value = 11207.4999
$ws.Range("M132").Value = -3635.1671
$ws.Range("N132").Value = -16267.4999
$ws.Range("H137").Value = 2564.8965
$ws.Range("I137").Value = 1483.3684
$ws.Range("J137").Value = 4619.8
$ws.Range("K137").Value = 4450.1052
$ws.Range("L137").Value = 13859.4
$ws.Range("M137").Value = -1900.1052
$ws.Range("N137").Value = -18959.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1309.875
$ws.Range("I74").Value = 1289.7097
$ws.Range("J74").Value = 1379.3334
$ws.Range("K74").Value = 1289.7097
$ws.Range("L74").Value = 1379.3334
$ws.Range("M74").Value = -415.7097000000001
$ws.Range("N74").Value = -3127.3334
$ws.Range("H77").Value = 1309.875
$ws.Range("I77").Value = 1289.7097
$ws.Range("J77").Value = 1379.3334
$ws.Range("K77").Value = 6448.548500000001
$ws.Range("L77").Value = 6896.666999999999
$ws.Range("M77").Value = -2080.548500000001
$ws.Range("N77").Value = -15632.667
$ws.Range("H122").Value = 1468.875
$ws.Range("I122").Value = 1359.5625
$ws.Range("J122").Value = 1687.5
$ws.Range("K122").Value = 4078.6875
$ws.Range("L122").Value = 5062.5
$ws.Range("M122").Value = -1628.6875
$ws.Range("N122").Value = -9962.5
$ws.Range("H132").Value = 2265.5715
$ws.Range("I132").Value = 1234.3077
$ws.Range("J132").Value = 3941.375
$ws.Range("K132").Value = 3702.9231
$ws.Range("L132").Value = 11824.125
$ws.Range("M132").Value = -1172.9231
$ws.Range("N132").Value = -16884.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1814.1111
$ws.Range("I31").Value = 1408.175
$ws.Range("J31").Value = 2973.9285
$ws.Range("K31").Value = 1408.175
$ws.Range("L31").Value = 2973.9285
$ws.Range("M31").Value = -1113.175
$ws.Range("N31").Value = -3563.9285
$ws.Range("H34").Value = 1814.1111
$ws.Range("I34").Value = 1408.175
$ws.Range("J34").Value = 2973.9285
$ws.Range("K34").Value = 1408.175
$ws.Range("L34").Value = 2973.9285
$ws.Range("M34").Value = -1206.175
$ws.Range("N34").Value = -3377.9285
$ws.Range("H58").Value = 814.6799999999999
$ws.Range("I58").Value = 628.0465
$ws.Range("J58").Value = 1961.1428
$ws.Range("K58").Value = 628.0465
$ws.Range("L58").Value = 1961.1428
$ws.Range("M58").Value = -425.0465
$ws.Range("N58").Value = -2367.1428
$ws.Range("H136").Value = 814.6799999999999
$ws.Range("I136").Value = 628.0465
$ws.Range("J136").Value = 1961.1428
$ws.Range("K136").Value = 1884.1395
$ws.Range("L136").Value = 5883.428400000001
$ws.Range("M136").Value = 665.8604999999998
$ws.Range("N136").Value = -10983.4284

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3216.5
$ws.Range("I75").Value = 300
$ws.Range("J75").Value = 3799.8
$ws.Range("K75").Value = 900
$ws.Range("L75").Value = 11399.4
$ws.Range("M75").Value = 98
$ws.Range("N75").Value = -13395.4
$ws.Range("H78").Value = 3216.5
$ws.Range("I78").Value = 300
$ws.Range("J78").Value = 3799.8
$ws.Range("K78").Value = 2700
$ws.Range("L78").Value = 34198.2
$ws.Range("M78").Value = 2292
$ws.Range("N78").Value = -44182.2
$ws.Range("H98").Value = 734.75
$ws.Range("J98").Value = 930.8
$ws.Range("L98").Value = 2792.4
$ws.Range("N98").Value = -5788.4
$ws.Range("H100").Value = 999
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H104").Value = 3399.889
$ws.Range("J104").Value = 3524.875
$ws.Range("L104").Value = 10574.625
$ws.Range("N104").Value = -15816.625
$ws.Range("H121").Value = 748.63635
$ws.Range("I121").Value = 441.875
$ws.Range("J121").Value = 1566.6666
$ws.Range("K121").Value = 1325.625
$ws.Range("L121").Value = 4699.9998
$ws.Range("M121").Value = -15.625
$ws.Range("N121").Value = -7319.9998
$ws.Range("H132").Value = 1444788.4
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1444788.4
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 13003095.6
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -13008155.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3131.2173
$ws.Range("I80").Value = 2884.6155
$ws.Range("J80").Value = 3451.8
$ws.Range("K80").Value = 2884.6155
$ws.Range("L80").Value = 3451.8
$ws.Range("M80").Value = -1886.6155
$ws.Range("N80").Value = -5447.8
$ws.Range("H83").Value = 3131.2173
$ws.Range("I83").Value = 2884.6155
$ws.Range("J83").Value = 3451.8
$ws.Range("K83").Value = 14423.0775
$ws.Range("L83").Value = 17259
$ws.Range("M83").Value = -9431.077499999999
$ws.Range("N83").Value = -27243
$ws.Range("H132").Value = 2158.0754
$ws.Range("I132").Value = 1887.7142
$ws.Range("J132").Value = 3190.3635
$ws.Range("K132").Value = 5663.142599999999
$ws.Range("L132").Value = 9571.0905
$ws.Range("M132").Value = -3133.142599999999
$ws.Range("N132").Value = -14631.0905

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 50707.5
$ws.Range("J133").Value = 50707.5
$ws.Range("L133").Value = 50707.5
$ws.Range("N133").Value = -55767.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2953.5293
$ws.Range("I122").Value = 2940.8333
$ws.Range("J122").Value = 2984
$ws.Range("K122").Value = 8822.499899999999
$ws.Range("L122").Value = 8952
$ws.Range("M122").Value = -6372.499899999999
$ws.Range("N122").Value = -13852
$ws.Range("H132").Value = 1110.0339
$ws.Range("I132").Value = 770
$ws.Range("J132").Value = 1950.1177
$ws.Range("K132").Value = 2310
$ws.Range("L132").Value = 5850.3531
$ws.Range("M132").Value = 220
$ws.Range("N132").Value = -10910.3531

Write-Host "Applied all Shinryu_Profits updates"